## Regression tool updated with printing summary results; tests added, docs updated
##
## Inserts a new "4.3. Paylod type inference" block of 11 test rows into the
## "RegressionTests" sheet, just above the existing "4.1.1. Tuples: duplicate
## names ..." section (which - along with everything below it - shifts down
## by 6 rows as a result).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RegressionTests")

# Make room: insert 6 new rows starting at row 49. Everything that used to
# live at row 49 onward (the "4.1.1. Tuples..." section header and below)
# moves down to row 55 onward.
$ws.Range("A49:A54").EntireRow.Insert()

# New payload-type-inference test rows (49-59).
$ws.Range("A49").Value = "payloadActions"
$ws.Range("B49").Value = "4.3. Paylod type inference"
$ws.Range("C49").Value = "Yes"

$ws.Range("A50").Value = "payloadActionsFuns"
$ws.Range("B50").Value = "4.3. Paylod type inference"
$ws.Range("C50").Value = "Yes"

$ws.Range("A51").Value = "payloadEntry"
$ws.Range("B51").Value = "4.3. Paylod type inference"
$ws.Range("C51").Value = "Yes"

$ws.Range("A52").Value = "payloadEntry_1"
$ws.Range("B52").Value = "4.3. Paylod type inference"
$ws.Range("C52").Value = "Yes"

$ws.Range("A53").Value = "payloadEntryFuns"
$ws.Range("B53").Value = "4.3. Paylod type inference"
$ws.Range("C53").Value = "Yes"

$ws.Range("A54").Value = "payloadExit"
$ws.Range("B54").Value = "4.3. Paylod type inference"
$ws.Range("C54").Value = "Yes"

$ws.Range("A55").Value = "payloadExitFuns"
$ws.Range("B55").Value = "4.3. Paylod type inference"
$ws.Range("C55").Value = "Yes"

$ws.Range("A56").Value = "payloads"
$ws.Range("B56").Value = "4.3. Paylod type inference"
$ws.Range("C56").Value = "Yes"
$ws.Range("G56").Value = '"invalid payload type in send"; "invalid payload type in send (cannot send null value)"; "invalid payload type in raise"; "argument 1 of "send" expects a machine value"; "argument 2 of "send" expects a machine value"'

$ws.Range("A57").Value = "payloadStartState"
$ws.Range("B57").Value = "4.3. Paylod type inference"
$ws.Range("C57").Value = "Yes"

$ws.Range("A58").Value = "payloadTransitions"
$ws.Range("B58").Value = "4.3. Paylod type inference"
$ws.Range("C58").Value = "Yes"

$ws.Range("A59").Value = "payloadTransitionsFuns"
$ws.Range("B59").Value = "4.3. Paylod type inference"
$ws.Range("C59").Value = "Yes"

# Match the updated view state recorded for the sheet after the edit.
$ws.Application.ActiveWindow.ScrollRow = 30
$ws.Range("G56").Select()
